# This script reproduces a sequence of Find & Replace edits that were made
# to the "Abstract" (D3) and "Authors" (E3) cells of the worksheet.
#
# D3 originally contained section markers left over from HTML scraping
# (id="Par1">, id="Par2">, id="Par3">, id="Par4">) that needed to be
# stripped out one at a time.
#
# E3 originally had its author list separators (",") repeatedly replaced
# with ", " (comma + space) -- an operation that (when run again on text
# that already contains a space after the comma) keeps adding one more
# space each time it's run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractCell = $ws.Cells.Item(3, 4)   # D3 - Abstract
$authorsCell  = $ws.Cells.Item(3, 5)   # E3 - Authors

# Strip the leftover id="ParN"> markers from the Abstract, one at a time.
$abstractCell.Replace('id="Par1">', '')
$abstractCell.Replace('id="Par2">', '')
$abstractCell.Replace('id="Par3">', '')
$abstractCell.Replace('id="Par4">', '')

# Re-run the comma -> comma+space normalization on the Authors list a few
# more times, which (since a space already follows each comma) pads the
# separators with additional spaces.
$authorsCell.Replace(',', ', ')
$authorsCell.Replace(',', ', ')
$authorsCell.Replace(',', ', ')
$authorsCell.Replace(',', ', ')
